# Update the generated LinkedIn draft slides with the new PFC Consulting
# article content (replacing the old Insolation Energy / Bondada Engineering
# article content) across all 6 carousel slides.

$p = $ppt.ActivePresentation

$newTitle = "PFC Consulting Accelerates Clean Energy Push with New 4.5 GW Transmission Project in Andhra Pradesh - SolarQuarter"
$newImageMissing = "⚠️ Image missing (not found in runner)"

$slideBodies = @{
    1 = @(
        "PFC Consulting has initiated a new transmission project with a capacity of 4.5 GW.",
        "The project is located in Andhra Pradesh."
    )
    2 = @(
        "The new transmission project aims to enhance clean energy infrastructure.",
        "It is part of PFC Consulting's broader strategy to support renewable energy initiatives."
    )
    3 = @(
        "The project is expected to facilitate the integration of renewable energy sources.",
        "It aligns with India's goals for increasing clean energy capacity."
    )
    4 = @(
        "The project contributes to the overall energy transition in the region.",
        "It is a significant step towards achieving energy sustainability."
    )
    5 = @(
        "PFC Consulting is focused on developing infrastructure for renewable energy.",
        "The 4.5 GW capacity will support various clean energy projects in Andhra Pradesh."
    )
    6 = @(
        "The project is part of a larger effort to enhance energy security in India.",
        "It reflects the growing investment in clean energy technologies."
    )
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Title placeholder (shape 1) - single paragraph
    $titleRange = $s.Shapes.Item(1).TextFrame.TextRange
    $titlePara = $titleRange.Paragraphs(1, 1)
    $titlePara.Text = ""
    $titleParaB = $titleRange.Paragraphs(1, 1)
    $titleParaB.Text = $newTitle

    # Content placeholder (shape 2) has 3 paragraphs:
    #   1) the "Image missing" warning line
    #   2) first bullet
    #   3) second bullet
    $bodyRange = $s.Shapes.Item(2).TextFrame.TextRange

    $para1 = $bodyRange.Paragraphs(1, 1)
    $para1.Text = ""
    $para1b = $bodyRange.Paragraphs(1, 1)
    $para1b.Text = $newImageMissing

    $bullets = $slideBodies[$i]

    $para2 = $bodyRange.Paragraphs(2, 1)
    $para2.Text = ""
    $para2b = $bodyRange.Paragraphs(2, 1)
    $para2b.Text = $bullets[0]

    $para3 = $bodyRange.Paragraphs(3, 1)
    $para3.Text = ""
    $para3b = $bodyRange.Paragraphs(3, 1)
    $para3b.Text = $bullets[1]
}
